$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price (column D) updates - keep values as text (leading apostrophe
# forces Excel to treat the numeric-looking string as text, matching the
# original inlineStr/text cell type in the workbook). Resetting the style
# back to Normal afterwards drops the auto-applied "quote prefix" style so
# the cell format stays identical to the original (unstyled) cells.
$updates = @{
    "D2"  = "264.88"
    "D3"  = "22.79"
    "D4"  = "6.223"
    "D5"  = "0.06170"
    "D6"  = "3.560"
    "D7"  = "6.699"
    "D8"  = "1.358"
    "D9"  = "0.8170"
    "D11" = "0.08208"
    "D13" = "0.03140"
    "D14" = "0.09254"
    "D15" = "3.888"
    "D16" = "0.001704"
    "D17" = "0.04844"
    "D18" = "0.0006256"
    "D19" = "0.006231"
    "D20" = "0.006276"
    "D23" = "3.702"
    "D24" = "2.261"
    "D25" = "0.3381"
    "D27" = "0.0002679"
    "D40" = "0.04592"
    "D44" = "0.01041"
    "D45" = "0.00006138"
    "D46" = "0.00000000749"
    "D47" = "0.7694"
    "D48" = "0.1952"
    "D49" = "0.00002098"
    "D50" = "0.01239"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}

# Rows 41-43 were reshuffled (coin order changed) with some values updated too.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1136"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003433"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICK"
